# Update TestCase Create/Edit/Delete service results:
# Mark BL01 and BL10-BL19 (rows 2, 10-20) as "PASS" in the Result column (E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$passRows = @(2, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 5).Value = "PASS"
}
